$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Prn"
$ws.Range("C2").Value = "Rpsa"
$ws.Range("G2").Value = 0.106052
$ws.Range("H2").Value = 0.318156
$ws.Range("I2").Value = 0.1062206075109891
$ws.Range("J2").Value = 0.1062206075109891
$ws.Range("M2").Value = 91.74689966666665
$ws.Range("N2").Value = 275.2406989999999
$ws.Range("O2").Value = 0.1908387282982634
$ws.Range("P2").Value = 0.1908387282982634
$ws.Range("Q2").Value = 9.729942203449331
$ws.Range("R2").Value = 87.56947983104398
$ws.Range("S2").Value = 0.02027100565646612
$ws.Range("T2").Value = 0.02027100565646613

# Row 3
$ws.Range("B3").Value = "Prn"
$ws.Range("C3").Value = "Rpsa"
$ws.Range("G3").Value = 0.106052
$ws.Range("H3").Value = 0.318156
$ws.Range("I3").Value = 0.1062206075109891
$ws.Range("J3").Value = 0.1062206075109891
$ws.Range("M3").Value = 142.7363333333333
$ws.Range("N3").Value = 428.209
$ws.Range("O3").Value = 0.296899627499751
$ws.Range("P3").Value = 0.296899627499751
$ws.Range("Q3").Value = 15.13747362266667
$ws.Range("R3").Value = 136.237262604
$ws.Range("S3").Value = 0.03153685880280991
$ws.Range("T3").Value = 0.03153685880280992

# Row 4
$ws.Range("B4").Value = "Prn"
$ws.Range("C4").Value = "Rpsa"
$ws.Range("G4").Value = 0.106052
$ws.Range("H4").Value = 0.318156
$ws.Range("I4").Value = 0.1062206075109891
$ws.Range("J4").Value = 0.1062206075109891
$ws.Range("M4").Value = 167.6324513333334
$ws.Range("N4").Value = 502.8973540000001
$ws.Range("O4").Value = 0.348684957750095
$ws.Range("P4").Value = 0.348684957750095
$ws.Range("Q4").Value = 17.77775672880267
$ws.Range("R4").Value = 159.999810559224
$ws.Range("S4").Value = 0.03703752804215866
$ws.Range("T4").Value = 0.03703752804215866

# Row 5
$ws.Range("B5").Value = "Prn"
$ws.Range("C5").Value = "Rpsa"
$ws.Range("G5").Value = 0.106052
$ws.Range("H5").Value = 0.318156
$ws.Range("I5").Value = 0.1062206075109891
$ws.Range("J5").Value = 0.1062206075109891
$ws.Range("M5").Value = 78.64050433333334
$ws.Range("N5").Value = 235.921513
$ws.Range("O5").Value = 0.1635766864518907
$ws.Range("P5").Value = 0.1635766864518907
$ws.Range("Q5").Value = 8.339982765558666
$ws.Range("R5").Value = 75.059844890028
$ws.Range("S5").Value = 0.01737521500955441
$ws.Range("T5").Value = 0.01737521500955441

# Row 6
$ws.Range("B6").Value = "Prn"
$ws.Range("C6").Value = "Rpsa"
$ws.Range("G6").Value = 0.7192416666666667
$ws.Range("H6").Value = 2.157725
$ws.Range("I6").Value = 0.7203851580408634
$ws.Range("J6").Value = 0.7203851580408636
$ws.Range("M6").Value = 91.74689966666665
$ws.Range("N6").Value = 275.2406989999999
$ws.Range("O6").Value = 0.1908387282982634
$ws.Range("P6").Value = 0.1908387282982634
$ws.Range("Q6").Value = 65.98819302775277
$ws.Range("R6").Value = 593.893737249775
$ws.Range("S6").Value = 0.1374773874454619
$ws.Range("T6").Value = 0.1374773874454619

# Row 7
$ws.Range("B7").Value = "Prn"
$ws.Range("C7").Value = "Rpsa"
$ws.Range("G7").Value = 0.7192416666666667
$ws.Range("H7").Value = 2.157725
$ws.Range("I7").Value = 0.7203851580408634
$ws.Range("J7").Value = 0.7203851580408636
$ws.Range("M7").Value = 142.7363333333333
$ws.Range("N7").Value = 428.209
$ws.Range("O7").Value = 0.296899627499751
$ws.Range("P7").Value = 0.296899627499751
$ws.Range("Q7").Value = 102.6619182805556
$ws.Range("R7").Value = 923.957264525
$ws.Range("S7").Value = 0.2138820850786816
$ws.Range("T7").Value = 0.2138820850786816

# Row 8
$ws.Range("B8").Value = "Prn"
$ws.Range("C8").Value = "Rpsa"
$ws.Range("G8").Value = 0.7192416666666667
$ws.Range("H8").Value = 2.157725
$ws.Range("I8").Value = 0.7203851580408634
$ws.Range("J8").Value = 0.7203851580408636
$ws.Range("M8").Value = 167.6324513333334
$ws.Range("N8").Value = 502.8973540000001
$ws.Range("O8").Value = 0.348684957750095
$ws.Range("P8").Value = 0.348684957750095
$ws.Range("Q8").Value = 120.5682436844056
$ws.Range("R8").Value = 1085.11419315965
$ws.Range("S8").Value = 0.251187468395274
$ws.Range("T8").Value = 0.251187468395274

# Row 9
$ws.Range("B9").Value = "Prn"
$ws.Range("C9").Value = "Rpsa"
$ws.Range("G9").Value = 0.7192416666666667
$ws.Range("H9").Value = 2.157725
$ws.Range("I9").Value = 0.7203851580408634
$ws.Range("J9").Value = 0.7203851580408636
$ws.Range("M9").Value = 78.64050433333334
$ws.Range("N9").Value = 235.921513
$ws.Range("O9").Value = 0.1635766864518907
$ws.Range("P9").Value = 0.1635766864518907
$ws.Range("Q9").Value = 56.56152740421389
$ws.Range("R9").Value = 509.053746637925
$ws.Range("S9").Value = 0.117838217121446
$ws.Range("T9").Value = 0.1178382171214461

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Prn"
$ws.Range("C10").Value = "Rpsa"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.173119
$ws.Range("H10").Value = 0.519357
$ws.Range("I10").Value = 0.1733942344481473
$ws.Range("J10").Value = 0.1733942344481473
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 91.74689966666665
$ws.Range("N10").Value = 275.2406989999999
$ws.Range("O10").Value = 0.1908387282982634
$ws.Range("P10").Value = 0.1908387282982634
$ws.Range("Q10").Value = 15.88313152339366
$ws.Range("R10").Value = 142.948183710543
$ws.Range("S10").Value = 0.03309033519633537
$ws.Range("T10").Value = 0.03309033519633537

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Prn"
$ws.Range("C11").Value = "Rpsa"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.173119
$ws.Range("H11").Value = 0.519357
$ws.Range("I11").Value = 0.1733942344481473
$ws.Range("J11").Value = 0.1733942344481473
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 142.7363333333333
$ws.Range("N11").Value = 428.209
$ws.Range("O11").Value = 0.296899627499751
$ws.Range("P11").Value = 0.296899627499751
$ws.Range("Q11").Value = 24.71037129033333
$ws.Range("R11").Value = 222.393341613
$ws.Range("S11").Value = 0.05148068361825943
$ws.Range("T11").Value = 0.05148068361825943

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Prn"
$ws.Range("C12").Value = "Rpsa"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.173119
$ws.Range("H12").Value = 0.519357
$ws.Range("I12").Value = 0.1733942344481473
$ws.Range("J12").Value = 0.1733942344481473
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 167.6324513333334
$ws.Range("N12").Value = 502.8973540000001
$ws.Range("O12").Value = 0.348684957750095
$ws.Range("P12").Value = 0.348684957750095
$ws.Range("Q12").Value = 29.02036234237534
$ws.Range("R12").Value = 261.183261081378
$ws.Range("S12").Value = 0.06045996131266232
$ws.Range("T12").Value = 0.06045996131266232

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Prn"
$ws.Range("C13").Value = "Rpsa"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.173119
$ws.Range("H13").Value = 0.519357
$ws.Range("I13").Value = 0.1733942344481473
$ws.Range("J13").Value = 0.1733942344481473
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 78.64050433333334
$ws.Range("N13").Value = 235.921513
$ws.Range("O13").Value = 0.1635766864518907
$ws.Range("P13").Value = 0.1635766864518907
$ws.Range("Q13").Value = 13.61416546968233
$ws.Range("R13").Value = 122.527489227141
$ws.Range("S13").Value = 0.02836325432089022
$ws.Range("T13").Value = 0.02836325432089022
